# Update header labels for area/subarea columns to include "Código"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Código Área (*)"
$ws.Range("F1").Value = "Código Subárea 1 (*)"
$ws.Range("G1").Value = "Código Subárea 2"
$ws.Range("H1").Value = "Código Subárea 3"

# Replace area/subarea descriptive names with their codes (row 2)
$ws.Range("E2").Value = "TEC"
$ws.Range("F2").Value = "COMP"
$ws.Range("H2").Value = "IEB"

# Replace area/subarea descriptive names with their codes (row 3)
$ws.Range("E3").Value = "BC"
$ws.Range("F3").Value = "BQM"
$ws.Range("H3").Value = "PAR"

# Update the active cell selection to F3
$ws.Range("F3").Select()
